$d = $word.ActiveDocument

# Locate the run containing "เชียงใหม่" and update its text, then insert
# a new run with the zip code "50300" right after it.
$found = $d.Content.Find.Execute("เชียงใหม่", $true, $false, $false, $false, $false,
                                  $true, 1, $false, "เชียงใหม่ 50300", 2)
